# Resize the table on slide 1 (shrink it and reposition to the origin).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

# Narrow the second column to match the first (1828800 EMU = 144 pt).
$tbl.Columns.Item(2).Width = 144

# Halve the row heights (2743200 EMU = 216 pt -> 1371600 EMU = 108 pt).
$tbl.Rows.Item(1).Height = 108
$tbl.Rows.Item(2).Height = 108

# Move the table's frame back to the top-left corner of the slide.
$sh.Left = 0
$sh.Top = 0
